$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, matching the formatting of the other header cells
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column values (0 for both data rows)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
